$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear "LM ING" (row 14, col A) and "Emmanuel Lalau" / "Marie-Anne Ferrand" entries
# that were removed from their original spot...
$ws.Range("A14").Value = ""
$ws.Range("B3").Value = ""
$ws.Range("B9").Value = ""

# ...Marie-Anne Ferrand is re-entered under the RH column (E7)
$ws.Range("E7").Value = "Marie-Anne Ferrand"

# New contact appended in the first free ENTREPRISE slot (row 17)
$ws.Range("A17").Value = "Julien Dugarry"

# Move the active selection to A5
$ws.Range("A5").Select() | Out-Null
